# The workbook gained a brand new review row. It was inserted as row 11
# (pushing the former rows 11-25 down to 12-26), so the sheet's used range
# grew from A1:D25 to A1:D26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 11, shifting existing rows 11-25 down
# to rows 12-26.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row with the new review data.
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = 45965.52277326389
$ws.Range("D11").Value = "MmJmMjU1YjgtYWExNC00YjFjLTg5NWQtMGYxNGVmNDA2ZDRhOjU3MDE2"
